# Update the LR-pairs sheet with newly computed TPM-based values.
# The underlying data now expands from 6 combinations (2 target clusters)
# to 9 combinations (3 target clusters: ECs, FAPs, MuSCs) for every
# sending cluster, and all numeric metrics are refreshed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Thbs1","Tnfrsf11b","ECs",3,1,2.727484333333333,8.182453000000001,0.03096049453772388,0.03096049453772388,1,0.3333333333333333,0.07580833333333332,0.227425,0.028190957994264,0.02819095799426401,0.2067660415027777,1.860894373525,0.0008728060009946141,0.0008728060009946143),
    @("ECs","Thbs1","Tnfrsf11b","FAPs",3,1,2.727484333333333,8.182453000000001,0.03096049453772388,0.03096049453772388,3,1,2.095195666666667,6.285587,0.7791435378093522,0.7791435378093522,5.714613356101222,51.431520204911,0.02412266924644931,0.02412266924644931),
    @("ECs","Thbs1","Tnfrsf11b","MuSCs",3,1,2.727484333333333,8.182453000000001,0.03096049453772388,0.03096049453772388,3,1,0.5180969999999999,1.554291,0.1926655041963838,0.1926655041963838,1.413101450647,12.717913055823,0.005965019290279958,0.00596501929027996),
    @("FAPs","Thbs1","Tnfrsf11b","ECs",3,1,52.73412466666667,158.202374,0.5986009007423507,0.5986009007423507,1,0.3333333333333333,0.07580833333333332,0.227425,0.028190957994264,0.02819095799426401,3.997686100772222,35.97917490694999,0.01687513284815621,0.01687513284815621),
    @("FAPs","Thbs1","Tnfrsf11b","FAPs",3,1,52.73412466666667,158.202374,0.5986009007423507,0.5986009007423507,3,1,2.095195666666667,6.285587,0.7791435378093522,0.7791435378093522,110.4883094870598,994.3947853835379,0.46639602354026,0.46639602354026),
    @("FAPs","Thbs1","Tnfrsf11b","MuSCs",3,1,52.73412466666667,158.202374,0.5986009007423507,0.5986009007423507,3,1,0.5180969999999999,1.554291,0.1926655041963838,0.1926655041963838,27.32139178742599,245.892526086834,0.1153297443539345,0.1153297443539345),
    @("MuSCs","Thbs1","Tnfrsf11b","ECs",3,1,32.63402300000001,97.90206900000001,0.3704386047199253,0.3704386047199253,1,0.3333333333333333,0.07580833333333332,0.227425,0.028190957994264,0.02819095799426401,2.473930893591667,22.265378042325,0.01044301914511318,0.01044301914511318),
    @("MuSCs","Thbs1","Tnfrsf11b","FAPs",3,1,32.63402300000001,97.90206900000001,0.3704386047199253,0.3704386047199253,3,1,2.095195666666667,6.285587,0.7791435378093522,0.7791435378093522,68.37466357550035,615.371972179503,0.2886248450226428,0.2886248450226428),
    @("MuSCs","Thbs1","Tnfrsf11b","MuSCs",3,1,32.63402300000001,97.90206900000001,0.3704386047199253,0.3704386047199253,3,1,0.5180969999999999,1.554291,0.1926655041963838,0.1926655041963838,16.907589414231,152.168304728079,0.07137074055216933,0.07137074055216934)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}

$ws.Range("A1:T10").Columns.AutoFit() | Out-Null
